$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.113.18"
$ws.Range("E2").Value = "  +2.99%  "
$ws.Range("D3").Value = "3.065.60"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.46"
$ws.Range("E5").Value = "  +2.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.09"
$ws.Range("E6").Value = "  +3.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.064.66"
$ws.Range("E8").Value = "  +1.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("E9").Value = "  +5.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("E10").Value = "  +5.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.14"
$ws.Range("E11").Value = "  -9.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.495"
$ws.Range("E12").Value = "  +11.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000232"
$ws.Range("E13").Value = "  +5.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.66"
$ws.Range("E14").Value = "  +5.04%  "
$ws.Range("D15").Value = "3.562.09"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").Value = "64.123.57"
$ws.Range("E16").Value = "  +3.16%  "
$ws.Range("D17").Value = "3.063.22"
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.80"
$ws.Range("E19").Value = "  +3.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.70"
$ws.Range("E20").Value = "  +2.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.01"
$ws.Range("E21").Value = "  +5.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.686"
$ws.Range("E22").Value = "  +4.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.51"
$ws.Range("E23").Value = "  +15.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.59"
$ws.Range("E24").Value = "  +4.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.58"
$ws.Range("E25").Value = "  +4.26%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.82"
$ws.Range("E27").Value = "  +4.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.20"
$ws.Range("E28").Value = "  +6.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.05"
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.38"
$ws.Range("E31").Value = "  +3.82%  "
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.46"
$ws.Range("E33").Value = "  +4.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.84"
$ws.Range("E34").Value = "  +5.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.26"
$ws.Range("E35").Value = "  +7.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.92"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0413"
$ws.Range("E37").Value = "  +5.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "451.46"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0819"
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.86"
$ws.Range("E40").Value = "  +11.31%  "
$ws.Range("D41").Value = "3.026.80"
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.31"
$ws.Range("E42").Value = "  +2.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.117"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.268"
$ws.Range("E44").Value = "  +8.84%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "28.00"
$ws.Range("E45").Value = "  +4.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.25"
$ws.Range("E46").Value = "  +13.66%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.114"
$ws.Range("E48").Value = "  +3.94%  "
$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").Value = "0.0₃0520"
$ws.Range("E49").Value = "  +3.44%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "118.36"
$ws.Range("E50").Value = "  +3.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.09"
$ws.Range("E51").Value = "  +4.91%  "
